$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D header and data
$ws.Range("D1").Value = "D"

# Fill D2:D4 values
$ws.Range("D2").Value = 0.2
$ws.Range("D3").Value = 0.2
$ws.Range("D4").Value = 1

# New row 4
$ws.Range("A4").Value = "D"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 5

# C2 should lose its number-format style (revert to default / General)
$ws.Range("C2").ClearFormats()

# Update selection to match target (C4)
$ws.Range("C4").Select()
